$d = $word.ActiveDocument

# 1. Title / H1 and the later bold "title" run (same text, 2 occurrences) -> replace all
$d.Content.Find.Execute("Play Leprechaun Legends Slot Game for Free | Pros & Cons Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Leprechaun Legends for free", 2)

# 2. "What we like" bullets - replace each distinct phrase (each unique in doc)
$d.Content.Find.Execute("Customizable betting range", $true, $false, $false, $false, $false, $true, 1, $false, "Stunning graphics set in an Irish forest", 2)
$d.Content.Find.Execute("Multiple bonus features", $true, $false, $false, $false, $false, $true, 1, $false, "Enchanting sound effects and background music", 2)
$d.Content.Find.Execute("Top-notch audiovisual experience", $true, $false, $false, $false, $false, $true, 1, $false, "Multiple bonus features for exciting gameplay", 2)
$d.Content.Find.Execute("Extra free games triggered by Leprechaun symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Customizable betting range for different player preferences", 2)

# 3. "What we don't like" bullet
$d.Content.Find.Execute("No progressive jackpot feature", $true, $false, $false, $false, $false, $true, 1, $false, "Betting range may not suit high rollers", 2)

# 4. Meta description (italic run)
$d.Content.Find.Execute("A detailed review of Leprechaun Legends, an Irish-themed online slot machine by Genesis with customizable betting range and multiple bonus features. Play for free!", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Leprechaun Legends and play this Irish-themed slot game for free.", 2)
